$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels -------------------------------------
# Columns A-J were "<Label>_old" and columns L-U were "<Label>_new".
# They become "<Label>_FV2310" and "<Label>_FV2404" respectively.
$labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $labels[$i] + "_FV2310"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $labels[$i] + "_FV2404"
}

# --- 2. Turn the data range into a native Excel Table ---------------------
$tableRange = $ws.Range("A1:U77")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
